$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 (shankejinjiboy) data values ---
$ws.Range("C5").Value = 224.41
$ws.Range("D5").Value = 163.98
$ws.Range("E5").Value = 15.44
$ws.Range("F5").Value = 18.75
$ws.Range("G5").Value = -0.05

# --- Apply new number formats ---
# Broad block (M_IBI, SD_IBI, MAE_HR, RMSE_HR, R_HR, MAE_RR, RMSE_RR, R_RR) gets
# the updated "no parentheses" negative format.
$ws.Range("C2:J10").NumberFormat = "0.00_ ;[Red]\-0.00\ "

# The HR/RR metrics block (RMSE_HR, R_HR, MAE_RR, RMSE_RR) switches to the plain
# two-decimal format without the red-negative styling.
$ws.Range("F3:I10").NumberFormat = "0.00_ "

# --- Update the active cell / selection cosmetics ---
$ws.Range("G8").Select()
